$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-12-12 Friday" "2025-12-13 Saturday"

Replace-Text "534÷5=" "378÷3="
Replace-Text "919÷7=" "791÷4="
Replace-Text "646÷4=" "731÷6="
Replace-Text "350÷5=" "899÷5="
Replace-Text "227÷2=" "348÷3="
Replace-Text "407÷8=" "584÷4="
Replace-Text "273÷9=" "617÷9="
Replace-Text "386÷6=" "960÷3="
Replace-Text "262÷7=" "389÷7="
Replace-Text "245÷6=" "904÷6="
Replace-Text "921÷8=" "469÷5="
Replace-Text "320÷8=" "751÷6="
Replace-Text "664÷7=" "853÷5="
Replace-Text "211÷4=" "296÷7="
Replace-Text "423÷4=" "133÷2="
Replace-Text "928÷8=" "626÷5="
Replace-Text "188÷8=" "469÷2="
Replace-Text "180÷2=" "204÷6="
Replace-Text "497÷6=" "586÷8="
Replace-Text "479÷8=" "974÷2="
Replace-Text "444÷7=" "757÷4="
Replace-Text "230÷4=" "117÷9="
Replace-Text "336÷6=" "607÷8="
Replace-Text "883÷6=" "360÷5="
Replace-Text "679÷8=" "923÷8="
